# The document has three inline logo pictures living in the headers/footers
# (the "first page" header/footer pair plus the "default" footer, since the
# section has Different First Page turned on):
#   - default footer   -> Pearson logo, currently named "image1.png"
#   - first-page footer-> Pearson logo, currently named "image1.png"
#   - first-page header-> BTec logo,    currently named "image2.jpg"
#
# The edit renames each picture:
#   Pearson logo (both footers): image1.png -> image2.png
#   BTec logo (header):          image2.jpg -> image1.jpg
#
# InlineShapes addressed straight off a Header/Footer Range can report a
# stale "addressed block" once we try to write to them here, so we first
# move the selection onto the shape's own Range and rename it through
# $word.Selection.InlineShapes - that resolves a live handle and the
# rename sticks.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-HeaderFooterPicture($range, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Default footer (wdHeaderFooterPrimary) - Pearson logo.
Rename-HeaderFooterPicture $sec.Footers.Item(1).Range "image2.png"

# First-page footer (wdHeaderFooterFirstPage) - Pearson logo.
Rename-HeaderFooterPicture $sec.Footers.Item(2).Range "image2.png"

# First-page header (wdHeaderFooterFirstPage) - BTec logo.
Rename-HeaderFooterPicture $sec.Headers.Item(2).Range "image1.jpg"
